# Dataset_HargaEmas_2025.xlsx -- "update manual book & dataset"
#
# Appends 7 new daily rows (2025-10-25 .. 2025-10-31) to the bottom of the
# Data_Harian_Lengkap sheet and moves the visible selection/scroll position
# to where the author was last working (around row 267), matching the
# window state recorded in the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows (dates are Excel serial numbers so they store as real
#     dates, the way Excel itself would write a freshly-typed date into a
#     cell that already carries the yyyy-mm-dd number format from column A).
#     Columns: Date(serial), Gold_Price, USD_Sell_Rate, USD_Buy_Rate, BI_Rate
$newRows = @(
    @(45955, 2350000, 16728, 16561, 0.0475),  # 2025-10-25
    @(45956, 2350000, 16728, 16561, 0.0475),  # 2025-10-26
    @(45957, 2327000, 16713, 16546, 0.0475),  # 2025-10-27
    @(45958, 2282000, 16711, 16544, 0.0475),  # 2025-10-28
    @(45959, 2267000, 16705, 16538, 0.0475),  # 2025-10-29
    @(45960, 2263000, 16714, 16547, 0.0475),  # 2025-10-30
    @(45961, 2305000, 16723, 16556, 0.0475)   # 2025-10-31
)

$startRow = 299
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

Write-Host "Appended $($newRows.Count) rows starting at row $startRow"

# --- Move the selection/scroll to where the workbook was left (B267, with
#     the viewport scrolled so row 247 is at the top) so the saved sheet
#     view matches the author's session.
$ws.Range("B267").Select()
$excel.ActiveWindow.ScrollRow = 247
$excel.ActiveWindow.ScrollColumn = 1

Write-Host "Selection now at $($excel.ActiveCell.Address())"
